$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CLIENTES list for the "RECURRENTE_400K-500K" group (row 9, column B):
# a new client id, "10425", was added to the end of the existing dot-separated list.
$ws.Range("B9").Value = "955.993.998.1001.1006.1009.10424.20103.20125.20310.20384.40151.50623.60159.60162.60192.60225.70103.70113.10425"

# Update the view: scroll so column B is the left-most visible column,
# and move the active selection to B10.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B10").Select()
